$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 194

# Update row 3 values (keep style, new values)
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 144

# Add new rows 4 and 5, copying the style used in column A (style index 1 -> border+bold+center)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 50

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 40

# Row 6 replaces what used to be row 3 (A=1, B=27 now instead of 11)
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 27

# Apply the same style as A2/A3 (style id 1) to the new A4:A6 cells
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4:A6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false
